# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2..63) holds a date stamp per forecast row that is supposed
# to represent the "as-of" / vintage date of the quarterly snapshot used to
# build that row's forecast. Those dates were incorrectly stamped with the
# 1st of the (correct) month; they should instead be stamped the 15th of the
# *following* month (i.e. shifted forward by one month, landing mid-month).
#
# This script walks every populated row in column A and re-derives the
# correct date from the existing (wrong) one, rewriting the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldSerial = $cell.Value2
    if ($oldSerial -eq $null) {
        continue
    }

    $oldDate = [DateTime]::FromOADate([double]$oldSerial)

    # Move to the 1st of the following month, then slide to the 15th.
    $nextMonthStart = $oldDate.AddMonths(1)
    $newDate = $nextMonthStart.AddDays(15 - $nextMonthStart.Day)

    $cell.Value = $newDate.ToOADate()
}
